$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 165, shifting existing rows 165:171 down to 166:172
$ws.Rows.Item(165).Insert()

# Populate the newly inserted row 165 with the new weekly price record
$ws.Range("A165").Value = 11
$ws.Range("B165").Value = "Vega Monumental Concepción"
$ws.Range("C165").Value = "Bíobío"
$ws.Range("D165").Value = 44706
$ws.Range("E165").Value = 8
$ws.Range("F165").Value = "Fruta"
$ws.Range("G165").Value = 100108
$ws.Range("H165").Value = "Tropicales y subtropicales"
$ws.Range("I165").Value = 100108005
$ws.Range("J165").Value = "Piña"
$ws.Range("K165").Value = "Caramelo"
$ws.Range("L165").Value = "Segunda"
$ws.Range("M165").Value = 200
$ws.Range("N165").Value = 16000
$ws.Range("O165").Value = 17000
$ws.Range("P165").Value = 16500
$ws.Range("Q165").Value = "$/caja 14 unidades"
$ws.Range("R165").Value = "Ecuador"
$ws.Range("S165").Value = 1179
$ws.Range("T165").Value = 14
